# paises.xlsx - "Update countries & provincias Spain"
#
# The source data got refreshed (7 Abr 2020, 19:52 -> 20:22). Because the sheet
# is kept sorted by "Casos totales" descending, two countries (Emiratos Arabes
# Unidos and San Martin (Parte Holandesa)) jumped past their neighbours, so the
# shared-string table was rebuilt in row order - shifting which country name a
# handful of rows show - while a wider set of rows simply got refreshed case
# counts. Apply both kinds of change cell-by-cell via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 20:22"

# Country-name shifts (rows whose country changed because of the re-sort)
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("A45").Value = "Finlandia"
$ws.Range("A46").Value = "Tailandia"
$ws.Range("A47").Value = "Panama"
$ws.Range("A146").Value = "San Martin (Parte Holandesa)"
$ws.Range("A147").Value = "Islas Caimanes"
$ws.Range("A148").Value = "Puerto Rico"
$ws.Range("A149").Value = "Zambia"
$ws.Range("A150").Value = "Bermudas"

# Refreshed case counts: Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes
# Row 4
$ws.Range("B4").Value = 386571
$ws.Range("C4").Value = 19567
$ws.Range("D4").Value = 21316
$ws.Range("E4").Value = 352981
$ws.Range("F4").Value = 9104
$ws.Range("G4").Value = 1403
$ws.Range("H4").Value = 12274
# Row 8
$ws.Range("B8").Value = 106739
$ws.Range("C8").Value = 3364
$ws.Range("D8").Value = 36081
$ws.Range("E8").Value = 68716
$ws.Range("F8").Value = 4895
$ws.Range("G8").Value = 132
$ws.Range("H8").Value = 1942
# Row 17
$ws.Range("B17").Value = 12616
$ws.Range("C17").Value = 319
$ws.Range("D17").Value = 4046
$ws.Range("E17").Value = 8327
$ws.Range("F17").Value = 243
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 243
# Row 25
$ws.Range("B25").Value = 5903
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 5782
$ws.Range("F25").Value = 78
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 89
# Row 44
$ws.Range("B44").Value = 2359
$ws.Range("C44").Value = 283
$ws.Range("D44").Value = 186
$ws.Range("E44").Value = 2161
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 12
# Row 45
$ws.Range("B45").Value = 2308
$ws.Range("C45").Value = 132
$ws.Range("D45").Value = 300
$ws.Range("E45").Value = 1974
$ws.Range("F45").Value = 81
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 34
# Row 46
$ws.Range("B46").Value = 2258
$ws.Range("C46").Value = 38
$ws.Range("D46").Value = 888
$ws.Range("E46").Value = 1343
$ws.Range("F46").Value = 61
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 27
# Row 47
$ws.Range("B47").Value = 2100
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 14
$ws.Range("E47").Value = 2031
$ws.Range("F47").Value = 88
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 55
# Row 76
$ws.Range("B76").Value = 697
$ws.Range("C76").Value = 35
$ws.Range("D76").Value = 51
$ws.Range("E76").Value = 640
$ws.Range("F76").Value = 16
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 6
# Row 110
$ws.Range("B110").Value = 196
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 46
$ws.Range("E110").Value = 147
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 3
# Row 146
$ws.Range("B146").Value = 40
$ws.Range("C146").Value = 3
$ws.Range("D146").Value = 1
$ws.Range("E146").Value = 33
$ws.Range("F146").Value = 2
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 6
# Row 147
$ws.Range("B147").Value = 39
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 1
$ws.Range("E147").Value = 37
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 1
# Row 148
$ws.Range("B148").Value = 39
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 36
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 2
# Row 149
$ws.Range("B149").Value = 39
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 7
$ws.Range("E149").Value = 31
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1
# Row 150
$ws.Range("B150").Value = 39
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 17
$ws.Range("E150").Value = 20
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 2
